$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new standard rows with their regex keys
$ws.Range("A23").Value = "UWC-1"
$ws.Range("D23").Value = "UWC\D*1"
$ws.Range("E23").Value = "Calib"
$ws.Range("F23").Value = "Carbonate"

$ws.Range("A24").Value = "UWW-1"
$ws.Range("D24").Value = "UWW\D*1"
$ws.Range("E24").Value = "Calib"
$ws.Range("F24").Value = "Something"

$ws.Range("F25").Select()
